$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 30 de Septiembre de 2020 a las 16:06"

# Update country statistics (new data pulled for this refresh)
$ws.Range("B4").Value = 7408066
$ws.Range("C4").Value = 1920
$ws.Range("D4").Value = 4649827
$ws.Range("E4").Value = 2547404
$ws.Range("G4").Value = 50
$ws.Range("H4").Value = 210835
$ws.Range("B5").Value = 6245404
$ws.Range("C5").Value = 21885
$ws.Range("D5").Value = 5206044
$ws.Range("E5").Value = 941599
$ws.Range("G5").Value = 232
$ws.Range("H5").Value = 97761
$ws.Range("B25").Value = 291191
$ws.Range("C25").Value = 725
$ws.Range("E25").Value = 25632
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 9559
$ws.Range("B52").Value = 75542
$ws.Range("C52").Value = 825
$ws.Range("D52").Value = 48530
$ws.Range("E52").Value = 25041
$ws.Range("G52").Value = 8
$ws.Range("H52").Value = 1971
$ws.Range("B68").Value = 40229
$ws.Range("C68").Value = 110
$ws.Range("D68").Value = 37954
$ws.Range("E68").Value = 1684
$ws.Range("G68").Value = 1
$ws.Range("H68").Value = 591
$ws.Range("B72").Value = 38529
$ws.Range("C72").Value = 151
$ws.Range("D72").Value = 24908
$ws.Range("E72").Value = 12910
$ws.Range("G72").Value = 4
$ws.Range("H72").Value = 711
$ws.Range("B76").Value = 33551
$ws.Range("C76").Value = 72
$ws.Range("E76").Value = 1266
$ws.Range("B79").Value = 27469
$ws.Range("C79").Value = 243
$ws.Range("D79").Value = 20616
$ws.Range("E79").Value = 5997
$ws.Range("G79").Value = 13
$ws.Range("H79").Value = 856
$ws.Range("B87").Value = 17977
$ws.Range("C87").Value = 191
$ws.Range("D87").Value = 14959
$ws.Range("E87").Value = 2279
$ws.Range("G87").Value = 2
$ws.Range("H87").Value = 739
$ws.Range("B92").Value = 14759
$ws.Range("C92").Value = 44
$ws.Range("D92").Value = 13959
$ws.Range("E92").Value = 468
$ws.Range("B93").Value = 13961
$ws.Range("C93").Value = 47
$ws.Range("E93").Value = 2497
$ws.Range("B96").Value = 13373
$ws.Range("C96").Value = 948
$ws.Range("D96").Value = 3755
$ws.Range("E96").Value = 9308
$ws.Range("G96").Value = 26
$ws.Range("H96").Value = 310

# Namibia/Malasia swap position (sorted by new totals) with refreshed figures
$ws.Range("A97").Value = "Namibia"
$ws.Range("B97").Value = 11265
$ws.Range("C97").Value = 125
$ws.Range("D97").Value = 9014
$ws.Range("E97").Value = 2130
$ws.Range("G97").Value = 0
$ws.Range("H97").Value = 121
$ws.Range("A98").Value = "Malasia"
$ws.Range("B98").Value = 11224
$ws.Range("C98").Value = 89
$ws.Range("D98").Value = 9967
$ws.Range("E98").Value = 1121
$ws.Range("G98").Value = 2
$ws.Range("H98").Value = 136

$ws.Range("B107").Value = 9769
$ws.Range("C107").Value = 43
$ws.Range("D107").Value = 8572
$ws.Range("E107").Value = 1121
$ws.Range("G107").Value = 1
$ws.Range("H107").Value = 76
$ws.Range("B142").Value = 3379
$ws.Range("C142").Value = 5
$ws.Range("E142").Value = 136

# Islandia/Sudan del Sur swap position with refreshed figures
$ws.Range("A148").Value = "Islandia"
$ws.Range("B148").Value = 2728
$ws.Range("C148").Value = 33
$ws.Range("D148").Value = 2167
$ws.Range("E148").Value = 551
$ws.Range("H148").Value = 10
$ws.Range("A149").Value = "Sudan del Sur"
$ws.Range("B149").Value = 2700
$ws.Range("D149").Value = 1290
$ws.Range("E149").Value = 1361
$ws.Range("H149").Value = 49

# Nueva Caledonia/Santa Lucia swap position (figures unchanged)
$ws.Range("A207").Value = "Nueva Caledonia"
$ws.Range("A208").Value = "Santa Lucia"
